$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VIN version bump: SYMBOL_2000 -> SYMBOL_2017 for every data row
$ws.Range("B2").Value = "SYMBOL_2017"
$ws.Range("B3").Value = "SYMBOL_2017"
$ws.Range("B4").Value = "SYMBOL_2017"
$ws.Range("B5").Value = "SYMBOL_2017"

# Replace the generic placeholder symbol codes (BI/PD/UM/MP_SYMBOL columns
# AC:AF) with row-specific codes.
$ws.Range("AC2").Value = "BI001"
$ws.Range("AD2").Value = "PD001"
$ws.Range("AE2").Value = "UM001"
$ws.Range("AF2").Value = "MP001"

$ws.Range("AC3").Value = "BI002"
$ws.Range("AD3").Value = "PD002"
$ws.Range("AE3").Value = "UM002"
$ws.Range("AF3").Value = "MP002"

$ws.Range("AC4").Value = "BI003"
$ws.Range("AD4").Value = "PD003"
$ws.Range("AE4").Value = "UM003"
$ws.Range("AF4").Value = "MP003"

$ws.Range("AC5").Value = "BI004"
$ws.Range("AD5").Value = "PD004"
$ws.Range("AE5").Value = "UM004"
$ws.Range("AF5").Value = "MP004"

# Move the view: scroll right so column X is the left-most visible column,
# then select AC2:AF5 (the block we just edited) as the active selection.
$excel.ActiveWindow.ScrollColumn = 24
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AC2:AF5").Select()
